$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.089.64'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '1.654.73'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.21'
$ws.Range("E5").Value = '  -0.49%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5276'
$ws.Range("E6").Value = '  +0.69%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2606'
$ws.Range("E8").Value = '  -3.87%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06335'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.44'
$ws.Range("E10").Value = '  -2.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07785'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.497'
$ws.Range("E12").Value = '  +1.02%  '
$ws.Range("D13").Value = '1.673.38'
$ws.Range("E13").Value = '  +1.28%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5482'
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").Value = '0.0₅8166'
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.40'
$ws.Range("E16").Value = '  +0.61%  '
$ws.Range("D17").Value = '26.110.26'
$ws.Range("E17").Value = '  -0.69%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.002'
$ws.Range("E18").Value = '  -0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.556'
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.58'
$ws.Range("E20").Value = '  -1.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.07'
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.030'
$ws.Range("E22").Value = '  -0.94%  '
$ws.Range("E23").Value = '  -0.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '141.27'
$ws.Range("E24").Value = '  +0.99%  '
$ws.Range("E25").Value = '  +0.52%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.270'
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.22'
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E28").Value = '  +1.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.05917'
$ws.Range("E29").Value = '  -3.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.279'
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.526'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.246'
$ws.Range("E32").Value = '  -1.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.569'
$ws.Range("E33").Value = '  -4.11%  '
$ws.Range("B34").Value = 'MXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.804'
$ws.Range("E34").Value = '  +0.48%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9488'
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.410'
$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5645'
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01611'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.819'
$ws.Range("E39").Value = '  -3.64%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.8467'
$ws.Range("E40").Value = '  -1.21%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '1.022.65'
$ws.Range("E42").Value = '  -0.15%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '102.03'
$ws.Range("E43").Value = '  +1.52%  '
$ws.Range("D44").Value = '1.799.19'
$ws.Range("E44").Value = '  -0.53%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '57.11'
$ws.Range("E45").Value = '  -1.62%  '
$ws.Range("B46").Value = 'Frax'
$ws.Range("C46").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.008'
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("B47").Value = 'Mantle'
$ws.Range("C47").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4288'
$ws.Range("E47").Value = '  +1.66%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.474'
$ws.Range("E48").Value = '  -1.07%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05151'
$ws.Range("E49").Value = '  -0.69%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.817'
$ws.Range("E50").Value = '  -3.56%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.09718'
$ws.Range("E51").Value = '  -0.91%  '
